# Update logic for saving result:
# The "Patient" worksheet gains a new "Barcode" column (K) next to the
# existing headers, and the active selection moves to K2 (the first data
# cell under the new header) ready for the next scanned value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patient")

$ws.Activate()

# New header cell - this also registers "Barcode" as a shared string.
$ws.Range("K1").Value = "Barcode"

# Move/save the active selection to the first cell under the new column.
$ws.Range("K2").Select() | Out-Null
